$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (incl. date number format / style) from row 3's styled cells
# so the new row re-uses the same style indices rather than minting new ones.
$ws.Range("A3").Copy($ws.Range("A4"))
$ws.Range("G3").Copy($ws.Range("G4"))

# Row 4 values
$ws.Range("A4").Value = 42633.679074074076
$ws.Range("B4").Value = $false
$ws.Range("C4").Value = 9974
$ws.Range("D4").Value = 10000
$ws.Range("E4").Value = 19.32
$ws.Range("F4").Value = 19.22
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = -0.52
$ws.Range("I4").Value = $false
